$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(15).Insert()

$ws.Range("A15").Value = 10
$ws.Range("B15").Value = "Vega Modelo de Temuco"
$ws.Range("C15").Value = "La Araucanía"
$ws.Range("D15").Value = 44764
$ws.Range("E15").Value = 9
$ws.Range("F15").Value = "Fruta"
$ws.Range("G15").Value = 100108
$ws.Range("H15").Value = "Tropicales y subtropicales"
$ws.Range("I15").Value = 100108003
$ws.Range("J15").Value = "Maracuyá"
$ws.Range("K15").Value = "Sin especificar"
$ws.Range("L15").Value = "Primera"
$ws.Range("M15").Value = 40
$ws.Range("N15").Value = 32000
$ws.Range("O15").Value = 34000
$ws.Range("P15").Value = 33000
$ws.Range("Q15").Value = '$/caja 18 kilos'
$ws.Range("R15").Value = "Región de Arica y Parinacota"
$ws.Range("S15").Value = 1833
$ws.Range("T15").Value = 18
